$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the old image-path columns (AD2:AF2) from the first data row of
# store_uat - these referenced shared strings that are being retired.
$ws1.Range("AD2:AF2").ClearContents()

# Add the new ar-uat theme rows (E3:E5) that replace the old single
# "img\storefront\ar-uat" value.
$ws1.Range("E3").Value = "ar-uat-starter"
$ws1.Range("E4").Value = "ar-uat-essential"
$ws1.Range("E5").Value = "ar-uat-esstial-orders"

# Update the saved view state for the sheet: scrolled over so column W is
# the left-most visible column, with AG15 selected.
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 23
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("AG15").Select()
